# Sheet gen_res in user inputs
#
# Adds a new "gen_res|pmax" worksheet (populated like the other generation /
# storage template tabs) positioned between "gen|pmax" and "storage|inflow",
# and registers it in the "ReadMe" overview table (new row, pushing the
# storage rows down).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Create the new "gen_res|pmax" data sheet and fill it with the same
#    layout used by the other component tabs (header row + 3 hourly rows).
# ---------------------------------------------------------------------
$newWs = $wb.Worksheets.Add()
$newWs.Name = "gen_res|pmax"

$newWs.Range("A1").Value = "Time\Id"
$newWs.Range("B1").Value = 1

$newWs.Range("A2").Value = 1
$newWs.Range("B2").Value = 1000

$newWs.Range("A3").Value = 2
$newWs.Range("B3").Value = 1000

$newWs.Range("A4").Value = 3
$newWs.Range("B4").Value = 1000

# Move it right after "gen|pmax" (i.e. right before "storage|inflow").
$newWs.Move($wb.Worksheets.Item("storage|inflow"))

# ---------------------------------------------------------------------
# 2. Register the new tab on the "ReadMe" summary sheet: insert a row
#    right above the (now shifted down) storage rows and fill it in the
#    same way as the existing "gen|pmax" row.
# ---------------------------------------------------------------------
$readme = $wb.Worksheets.Item("ReadMe")
$readme.Rows.Item(10).Insert()

$readme.Range("A10").Value = "gen_res|pmax"
$readme.Range("B10").Formula = "=COUNT('gen_res|pmax'!`$1:`$1)"
$readme.Range("C10").Formula = "=IF(B10,COUNT('gen_res|pmax'!B:B)-1,0)"
$readme.Range("D10").Value = "MW"
$readme.Range("E10").Value = "Available production. If MW: in MW.`nIf %: 0.5 means 50% of the production rating."
$readme.Rows.Item(10).RowHeight = 30

# Restore the selection on the ReadMe sheet to match the saved view.
$readme.Activate() | Out-Null
$readme.Range("B10").Select() | Out-Null
